$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 5412.6924
$ws.Range("J88").Value = 5574.1113
$ws.Range("L88").Value = 5574.1113
$ws.Range("N88").Value = -6386.1113
$ws.Range("H91").Value = 5412.6924
$ws.Range("J91").Value = 5574.1113
$ws.Range("L91").Value = 5574.1113
$ws.Range("N91").Value = -8382.1113
$ws.Range("H100").Value = 2099.75
$ws.Range("I100").Value = 1775
$ws.Range("K100").Value = 1775
$ws.Range("M100").Value = -1234
$ws.Range("H101").Value = 717.1429000000001
$ws.Range("I101").Value = 246
$ws.Range("J101").Value = 1345.3334
$ws.Range("K101").Value = 738
$ws.Range("L101").Value = 4036.0002
$ws.Range("M101").Value = 884
$ws.Range("N101").Value = -7280.0002
$ws.Range("H113").Value = 4051.6924
$ws.Range("I113").Value = 3549.0527
$ws.Range("K113").Value = 3549.0527
$ws.Range("M113").Value = -295.0527000000002
$ws.Range("H138").Value = 7426.725
$ws.Range("I138").Value = 1483.0834
$ws.Range("K138").Value = 4449.2502
$ws.Range("M138").Value = 690.7497999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 83335500
$ws.Range("I61").Value = 100001800
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 100001800
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -100001588
$ws.Range("N61").Value = -4424
$ws.Range("H102").Value = 2853.3157
$ws.Range("I102").Value = 2789.611
$ws.Range("K102").Value = 2789.611
$ws.Range("M102").Value = -1167.611
$ws.Range("H136").Value = 83335500
$ws.Range("I136").Value = 100001800
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 300005400
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -300002850
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 40428
$ws.Range("I82").Value = 20000
$ws.Range("J82").Value = 50642
$ws.Range("K82").Value = 20000
$ws.Range("L82").Value = 50642
$ws.Range("M82").Value = -19617
$ws.Range("N82").Value = -51408
$ws.Range("H85").Value = 40428
$ws.Range("I85").Value = 20000
$ws.Range("J85").Value = 50642
$ws.Range("K85").Value = 20000
$ws.Range("L85").Value = 50642
$ws.Range("M85").Value = -18674
$ws.Range("N85").Value = -53294
$ws.Range("H94").Value = 2198.625
$ws.Range("I94").Value = 1181.8334
$ws.Range("K94").Value = 1181.8334
$ws.Range("M94").Value = -730.8334
$ws.Range("H95").Value = 40299.2
$ws.Range("J95").Value = 40299.2
$ws.Range("L95").Value = 40299.2
$ws.Range("N95").Value = -45791.2
$ws.Range("H99").Value = 2354.3333
$ws.Range("I99").Value = 2050.35
$ws.Range("J99").Value = 3874.25
$ws.Range("K99").Value = 2050.35
$ws.Range("L99").Value = 3874.25
$ws.Range("M99").Value = -552.3499999999999
$ws.Range("N99").Value = -6870.25
$ws.Range("H134").Value = 2203.1082
$ws.Range("I134").Value = 2283.4
$ws.Range("J134").Value = 798
$ws.Range("K134").Value = 6850.200000000001
$ws.Range("L134").Value = 2394
$ws.Range("M134").Value = -4315.200000000001
$ws.Range("N134").Value = -7464

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 35236.5
$ws.Range("J28").Value = 35236.5
$ws.Range("L28").Value = 35236.5
$ws.Range("N28").Value = -35726.5
$ws.Range("H58").Value = 5536.273
$ws.Range("I58").Value = 5489.9
$ws.Range("J58").Value = 6000
$ws.Range("K58").Value = 5489.9
$ws.Range("L58").Value = 6000
$ws.Range("M58").Value = -5286.9
$ws.Range("N58").Value = -6406
$ws.Range("H62").Value = 3172.389
$ws.Range("I62").Value = 3203.1177
$ws.Range("J62").Value = 2650
$ws.Range("K62").Value = 3203.1177
$ws.Range("L62").Value = 2650
$ws.Range("M62").Value = -2579.1177
$ws.Range("N62").Value = -3898
$ws.Range("H65").Value = 3172.389
$ws.Range("I65").Value = 3203.1177
$ws.Range("J65").Value = 2650
$ws.Range("K65").Value = 16015.5885
$ws.Range("L65").Value = 13250
$ws.Range("M65").Value = -12895.5885
$ws.Range("N65").Value = -19490
$ws.Range("H136").Value = 5536.273
$ws.Range("I136").Value = 5489.9
$ws.Range("J136").Value = 6000
$ws.Range("K136").Value = 16469.7
$ws.Range("L136").Value = 18000
$ws.Range("M136").Value = -13919.7
$ws.Range("N136").Value = -23100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 71.07692
$ws.Range("I2").Value = 28.333334
$ws.Range("J2").Value = 167.25
$ws.Range("K2").Value = 170.000004
$ws.Range("L2").Value = 1003.5
$ws.Range("M2").Value = -57.00000399999999
$ws.Range("N2").Value = -1229.5
$ws.Range("H7").Value = 1721.5
$ws.Range("I7").Value = 65.8
$ws.Range("K7").Value = 197.4
$ws.Range("M7").Value = -85.39999999999998
$ws.Range("H92").Value = 799.9091
$ws.Range("I92").Value = 190.14285
$ws.Range("J92").Value = 1867
$ws.Range("K92").Value = 570.4285500000001
$ws.Range("L92").Value = 5601
$ws.Range("M92").Value = 677.5714499999999
$ws.Range("N92").Value = -8097

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4581.552
$ws.Range("I70").Value = 4047.0833
$ws.Range("K70").Value = 4047.0833
$ws.Range("M70").Value = -3777.0833
$ws.Range("H73").Value = 4581.552
$ws.Range("I73").Value = 4047.0833
$ws.Range("K73").Value = 4047.0833
$ws.Range("M73").Value = -3111.0833
$ws.Range("H80").Value = 5691.3125
$ws.Range("I80").Value = 6219.222
$ws.Range("K80").Value = 6219.222
$ws.Range("M80").Value = -5221.222
$ws.Range("H83").Value = 5691.3125
$ws.Range("I83").Value = 6219.222
$ws.Range("K83").Value = 31096.11
$ws.Range("M83").Value = -26104.11
$ws.Range("H113").Value = 2703.7144
$ws.Range("I113").Value = 2079.2173
$ws.Range("J113").Value = 3459.6843
$ws.Range("K113").Value = 2079.2173
$ws.Range("L113").Value = 3459.6843
$ws.Range("M113").Value = 90.7827000000002
$ws.Range("N113").Value = -7799.6843
$ws.Range("H122").Value = 41668708
$ws.Range("I122").Value = 2155.125
$ws.Range("K122").Value = 6465.375
$ws.Range("M122").Value = -4015.375
$ws.Range("H132").Value = 3256.8
$ws.Range("I132").Value = 3325.8235
$ws.Range("K132").Value = 9977.470499999999
$ws.Range("M132").Value = -7447.470499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2733.2104
$ws.Range("J40").Value = 5016
$ws.Range("L40").Value = 5016
$ws.Range("N40").Value = -5288
$ws.Range("H122").Value = 4468645.5
$ws.Range("I122").Value = 4088.9
$ws.Range("K122").Value = 12266.7
$ws.Range("M122").Value = -9816.700000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 6.5
$ws.Range("I8").Value = 6.5
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 6.5
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 133.5
$ws.Range("N8").ClearContents()
$ws.Range("H100").Value = 400698.3
$ws.Range("J100").Value = 611.625
$ws.Range("L100").Value = 1223.25
$ws.Range("N100").Value = -2305.25
$ws.Range("H107").Value = 1596.32
$ws.Range("I107").Value = 714.6923
$ws.Range("J107").Value = 2551.4167
$ws.Range("K107").Value = 2144.0769
$ws.Range("L107").Value = 7654.250100000001
$ws.Range("M107").Value = -224.0769
$ws.Range("N107").Value = -11494.2501
$ws.Range("H109").Value = 92249.75
$ws.Range("I109").Value = 89666.336
$ws.Range("J109").Value = 100000
$ws.Range("K109").Value = 89666.336
$ws.Range("L109").Value = 100000
$ws.Range("M109").Value = -88279.336
$ws.Range("N109").Value = -102774
$ws.Range("H135").Value = 25035248
$ws.Range("J135").Value = 25035248
$ws.Range("L135").Value = 25035248
$ws.Range("N135").Value = -25045388
